# Move the first "available" name from Sheet1 (names) to the "used" sheet,
# recording which source image it was used for and when.

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# The name that is being consumed (row 2 of the names sheet, right after
# the header-like first entry).
$consumedName = $namesSheet.Range("A2").Value()

# Remove that name from the names list — everything below shifts up one row,
# shrinking the sheet's used range from A1:A475 to A1:A474.
$namesSheet.Rows.Item(2).Delete()

# Append a new record to the "used" tracking sheet describing this name's use.
$lastRow = $usedSheet.UsedRange.Rows.Count
$newRow = $lastRow + 1

$usedSheet.Range("A" + $newRow).Value = $consumedName
$usedSheet.Range("B" + $newRow).Value = "ChatGPT Image 2026年1月21日 14_13_20.png"
$usedSheet.Range("C" + $newRow).Value = "2026-01-21 14:15:31"
